$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated NATMI metrics (new TPM input) for rows 2-10, columns E:T
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.007258333333333333
$ws.Range("H2").Value = 0.021775
$ws.Range("I2").Value = 0.000328667160253549
$ws.Range("J2").Value = 0.000328667160253549
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007258333333333333
$ws.Range("N2").Value = 0.021775
$ws.Range("O2").Value = 0.000328667160253549
$ws.Range("P2").Value = 0.000328667160253549
$ws.Range("Q2").Value = 0.00005268340277777777
$ws.Range("R2").Value = 0.0004741506249999999
$ws.Range("S2").Value = 0.0000001080221022291321
$ws.Range("T2").Value = 0.000000108022102229132
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.007258333333333333
$ws.Range("H3").Value = 0.021775
$ws.Range("I3").Value = 0.000328667160253549
$ws.Range("J3").Value = 0.000328667160253549
$ws.Range("O3").Value = 0.7778551418094273
$ws.Range("P3").Value = 0.7778551418094272
$ws.Range("Q3").Value = 0.1246855807166667
$ws.Range("R3").Value = 1.12217022645
$ws.Range("S3").Value = 0.0002556554405471262
$ws.Range("T3").Value = 0.000255655440547126
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.007258333333333333
$ws.Range("H4").Value = 0.021775
$ws.Range("I4").Value = 0.000328667160253549
$ws.Range("J4").Value = 0.000328667160253549
$ws.Range("M4").Value = 4.898620999999999
$ws.Range("N4").Value = 14.695863
$ws.Range("O4").Value = 0.2218161910303192
$ws.Range("P4").Value = 0.2218161910303192
$ws.Range("Q4").Value = 0.03555582409166666
$ws.Range("R4").Value = 0.320002416825
$ws.Range("S4").Value = 0.00007290369760419377
$ws.Range("T4").Value = 0.00007290369760419376
$ws.Range("I5").Value = 0.7778551418094273
$ws.Range("J5").Value = 0.7778551418094272
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.007258333333333333
$ws.Range("N5").Value = 0.021775
$ws.Range("O5").Value = 0.000328667160253549
$ws.Range("P5").Value = 0.000328667160253549
$ws.Range("Q5").Value = 0.1246855807166667
$ws.Range("R5").Value = 1.12217022645
$ws.Range("S5").Value = 0.0002556554405471262
$ws.Range("T5").Value = 0.000255655440547126
$ws.Range("I6").Value = 0.7778551418094273
$ws.Range("J6").Value = 0.7778551418094272
$ws.Range("O6").Value = 0.7778551418094273
$ws.Range("P6").Value = 0.7778551418094272
$ws.Range("S6").Value = 0.6050586216393642
$ws.Range("T6").Value = 0.6050586216393641
$ws.Range("I7").Value = 0.7778551418094273
$ws.Range("J7").Value = 0.7778551418094272
$ws.Range("M7").Value = 4.898620999999999
$ws.Range("N7").Value = 14.695863
$ws.Range("O7").Value = 0.2218161910303192
$ws.Range("P7").Value = 0.2218161910303192
$ws.Range("Q7").Value = 84.149814571186
$ws.Range("R7").Value = 757.348331140674
$ws.Range("S7").Value = 0.172540864729516
$ws.Range("T7").Value = 0.1725408647295159
$ws.Range("G8").Value = 4.898620999999999
$ws.Range("H8").Value = 14.695863
$ws.Range("I8").Value = 0.2218161910303192
$ws.Range("J8").Value = 0.2218161910303192
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.007258333333333333
$ws.Range("N8").Value = 0.021775
$ws.Range("O8").Value = 0.000328667160253549
$ws.Range("P8").Value = 0.000328667160253549
$ws.Range("Q8").Value = 0.03555582409166666
$ws.Range("R8").Value = 0.320002416825
$ws.Range("S8").Value = 0.00007290369760419377
$ws.Range("T8").Value = 0.00007290369760419376
$ws.Range("G9").Value = 4.898620999999999
$ws.Range("H9").Value = 14.695863
$ws.Range("I9").Value = 0.2218161910303192
$ws.Range("J9").Value = 0.2218161910303192
$ws.Range("O9").Value = 0.7778551418094273
$ws.Range("P9").Value = 0.7778551418094272
$ws.Range("Q9").Value = 84.149814571186
$ws.Range("R9").Value = 757.348331140674
$ws.Range("S9").Value = 0.172540864729516
$ws.Range("T9").Value = 0.1725408647295159
$ws.Range("G10").Value = 4.898620999999999
$ws.Range("H10").Value = 14.695863
$ws.Range("I10").Value = 0.2218161910303192
$ws.Range("J10").Value = 0.2218161910303192
$ws.Range("M10").Value = 4.898620999999999
$ws.Range("N10").Value = 14.695863
$ws.Range("O10").Value = 0.2218161910303192
$ws.Range("P10").Value = 0.2218161910303192
$ws.Range("Q10").Value = 23.99648770164099
$ws.Range("R10").Value = 215.968389314769
$ws.Range("S10").Value = 0.04920242260319908
$ws.Range("T10").Value = 0.04920242260319907
